$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("values")

# Fix shared string label "M€/kTCO2" -> "M€/ktCO2" (used by B10)
$ws.Range("B10").Value = "M€/ktCO2"

# Reindexed CO2 budget values
$ws.Range("B3").Value = 0.03
$ws.Range("B9").Value = 0.04

# Update selection to B8
$ws.Range("B8").Select()
